# Update NATMI TPM-derived metrics (ligand/receptor expression & specificity
# columns E:T) on rows 2-10 of Sheet1 to reflect the refreshed TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2564746666666666
$ws.Cells.Item(2, 8).Value = 0.7694239999999999
$ws.Cells.Item(2, 9).Value = 0.1818007399394835
$ws.Cells.Item(2, 10).Value = 0.1818007399394835
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.05968133333333333
$ws.Cells.Item(2, 14).Value = 0.179044
$ws.Cells.Item(2, 15).Value = 0.02602747651633847
$ws.Cells.Item(2, 16).Value = 0.02602747651633848
$ws.Cells.Item(2, 17).Value = 0.01530675007288889
$ws.Cells.Item(2, 18).Value = 0.137760750656
$ws.Cells.Item(2, 19).Value = 0.004731814489427865
$ws.Cells.Item(2, 20).Value = 0.004731814489427866
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2564746666666666
$ws.Cells.Item(3, 8).Value = 0.7694239999999999
$ws.Cells.Item(3, 9).Value = 0.1818007399394835
$ws.Cells.Item(3, 10).Value = 0.1818007399394835
$ws.Cells.Item(3, 15).Value = 0.144012433133819
$ws.Cells.Item(3, 16).Value = 0.144012433133819
$ws.Cells.Item(3, 17).Value = 0.08469366286755554
$ws.Cells.Item(3, 18).Value = 0.7622429658079999
$ws.Cells.Item(3, 19).Value = 0.02618156690421369
$ws.Cells.Item(3, 20).Value = 0.02618156690421369
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2564746666666666
$ws.Cells.Item(4, 8).Value = 0.7694239999999999
$ws.Cells.Item(4, 9).Value = 0.1818007399394835
$ws.Cells.Item(4, 10).Value = 0.1818007399394835
$ws.Cells.Item(4, 15).Value = 0.8299600903498424
$ws.Cells.Item(4, 16).Value = 0.8299600903498425
$ws.Cells.Item(4, 17).Value = 0.4880992464053333
$ws.Cells.Item(4, 18).Value = 4.392893217647999
$ws.Cells.Item(4, 19).Value = 0.1508873585458419
$ws.Cells.Item(4, 20).Value = 0.150887358545842
$ws.Cells.Item(5, 9).Value = 0.7694380609030022
$ws.Cells.Item(5, 10).Value = 0.7694380609030022
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.05968133333333333
$ws.Cells.Item(5, 14).Value = 0.179044
$ws.Cells.Item(5, 15).Value = 0.02602747651633847
$ws.Cells.Item(5, 16).Value = 0.02602747651633848
$ws.Cells.Item(5, 17).Value = 0.06478299317555555
$ws.Cells.Item(5, 18).Value = 0.5830469385799999
$ws.Cells.Item(5, 19).Value = 0.0200265310609299
$ws.Cells.Item(5, 20).Value = 0.0200265310609299
$ws.Cells.Item(6, 9).Value = 0.7694380609030022
$ws.Cells.Item(6, 10).Value = 0.7694380609030022
$ws.Cells.Item(6, 15).Value = 0.144012433133819
$ws.Cells.Item(6, 16).Value = 0.144012433133819
$ws.Cells.Item(6, 19).Value = 0.1108086472964089
$ws.Cells.Item(6, 20).Value = 0.1108086472964089
$ws.Cells.Item(7, 9).Value = 0.7694380609030022
$ws.Cells.Item(7, 10).Value = 0.7694380609030022
$ws.Cells.Item(7, 15).Value = 0.8299600903498424
$ws.Cells.Item(7, 16).Value = 0.8299600903498425
$ws.Cells.Item(7, 19).Value = 0.6386028825456632
$ws.Cells.Item(7, 20).Value = 0.6386028825456633
$ws.Cells.Item(8, 7).Value = 0.06878966666666667
$ws.Cells.Item(8, 9).Value = 0.0487611991575143
$ws.Cells.Item(8, 10).Value = 0.0487611991575143
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 0.05968133333333333
$ws.Cells.Item(8, 14).Value = 0.179044
$ws.Cells.Item(8, 15).Value = 0.02602747651633847
$ws.Cells.Item(8, 16).Value = 0.02602747651633848
$ws.Cells.Item(8, 17).Value = 0.004105459026222222
$ws.Cells.Item(8, 18).Value = 0.03694913123599999
$ws.Cells.Item(8, 19).Value = 0.001269130965980707
$ws.Cells.Item(8, 20).Value = 0.001269130965980707
$ws.Cells.Item(9, 7).Value = 0.06878966666666667
$ws.Cells.Item(9, 9).Value = 0.0487611991575143
$ws.Cells.Item(9, 10).Value = 0.0487611991575143
$ws.Cells.Item(9, 15).Value = 0.144012433133819
$ws.Cells.Item(9, 16).Value = 0.144012433133819
$ws.Cells.Item(9, 17).Value = 0.02271588423588889
$ws.Cells.Item(9, 19).Value = 0.00702221893319636
$ws.Cells.Item(9, 20).Value = 0.00702221893319636
$ws.Cells.Item(10, 7).Value = 0.06878966666666667
$ws.Cells.Item(10, 9).Value = 0.0487611991575143
$ws.Cells.Item(10, 10).Value = 0.0487611991575143
$ws.Cells.Item(10, 15).Value = 0.8299600903498424
$ws.Cells.Item(10, 16).Value = 0.8299600903498425
$ws.Cells.Item(10, 19).Value = 0.04046984925833723
$ws.Cells.Item(10, 20).Value = 0.04046984925833724

Write-Output "done"
